$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits between
#    "controller" and " Systems" in the course-title line.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) "a.)" question: insert "(CLI) " in front of "development
#    environment", re-splitting the run in the process, and put the
#    "_GoBack" bookmark right before the newly inserted text.
#
# We force run splits by dropping a throw-away bookmark at the desired
# offset and immediately deleting it again: Word (and this host) always
# breaks the run around a bookmark, and the break persists even once
# the bookmark itself is removed.
# ---------------------------------------------------------------------
$full = $d.Content.Text
$needleA = "a.)  Describe the advantages of a command-line interface "
$idxA = $full.IndexOf($needleA)
$posA = $idxA + $needleA.Length

$rA = $d.Range($posA, $posA)
$rA.InsertAfter("(CLI) ")

# Split "(CLI) development environment.  " into "(CLI)" | " " | "development environment.  "
$splitA1 = $d.Range(($posA + 5), ($posA + 5))
$d.Bookmarks.Add("zzsplitA1", $splitA1)
$d.Bookmarks("zzsplitA1").Delete()

$splitA2 = $d.Range(($posA + 6), ($posA + 6))
$d.Bookmarks.Add("zzsplitA2", $splitA2)
$d.Bookmarks("zzsplitA2").Delete()

# Re-create the "_GoBack" bookmark right before "(CLI)"
$bmRangeA = $d.Range($posA, $posA)
$d.Bookmarks.Add("_GoBack", $bmRangeA)

# ---------------------------------------------------------------------
# 3) "b.)" question: insert " (IDE" right after "environment" (before
#    the trailing ".  "), splitting that run into three pieces, and add
#    a left tab stop at 7200 twips (360 pt) to the paragraph.
# ---------------------------------------------------------------------
$full = $d.Content.Text
$needleB = "b.)  Describe the advantages of an integrated development environment"
$idxB = $full.IndexOf($needleB)
$posB = $idxB + $needleB.Length

$rB = $d.Range($posB, $posB)
$rB.InsertAfter(" (IDE")

# Split "environment" | " (IDE" | ".  "
$splitB1 = $d.Range($posB, $posB)
$d.Bookmarks.Add("zzsplitB1", $splitB1)
$d.Bookmarks("zzsplitB1").Delete()

$splitB2 = $d.Range(($posB + 5), ($posB + 5))
$d.Bookmarks.Add("zzsplitB2", $splitB2)
$d.Bookmarks("zzsplitB2").Delete()

# Add the paragraph tab stop (7200 twips == 360 points).
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*integrated development environment*") {
        $p.Format.TabStops.Add(360)
        break
    }
}
